$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 17:24"

# --- Update country names (column A) that moved rows due to re-sorting ---
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("A46").Value = "Irlanda"
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("A81").Value = "Haiti"
$ws.Range("A82").Value = "Guinea"
$ws.Range("A86").Value = "Etiopia"
$ws.Range("A87").Value = "Gabon"
$ws.Range("A88").Value = "Luxemburgo"
$ws.Range("A89").Value = "Hungria"
$ws.Range("A102").Value = "Guayana Francesa"
$ws.Range("A103").Value = "Maldivas"
$ws.Range("A104").Value = "Costa Rica"
$ws.Range("A105").Value = "Estonia"
$ws.Range("A133").Value = "Cabo Verde"
$ws.Range("A134").Value = "Republica del Chad"
$ws.Range("A135").Value = "Principado de Andorra"
$ws.Range("A136").Value = "Uruguay"
$ws.Range("A210").Value = "Montserrat"
$ws.Range("A211").Value = "Seychelles"
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- Update statistic columns (B,C,D,E,F,G,H) with refreshed data ---
$ws.Range("B4").Value = 2300255 ; $ws.Range("C4").Value = 3065 ; $ws.Range("D4").Value = 956080 ; $ws.Range("E4").Value = 1222707 ; $ws.Range("G4").Value = 61 ; $ws.Range("H4").Value = 121468
$ws.Range("B7").Value = 400724 ; $ws.Range("C7").Value = 4912 ; $ws.Range("D7").Value = 216730 ; $ws.Range("E7").Value = 170959 ; $ws.Range("G7").Value = 65 ; $ws.Range("H7").Value = 13035
$ws.Range("B14").Value = 190798 ; $ws.Range("C14").Value = 138 ; $ws.Range("E14").Value = 7438
$ws.Range("D34").Value = 34224 ; $ws.Range("E34").Value = 7583
$ws.Range("B45").Value = 25778 ; $ws.Range("C45").Value = 710 ; $ws.Range("D45").Value = 14957 ; $ws.Range("E45").Value = 10166 ; $ws.Range("G45").Value = 8 ; $ws.Range("H45").Value = 655
$ws.Range("B46").Value = 25368 ; $ws.Range("D46").Value = 22698 ; $ws.Range("E46").Value = 956 ; $ws.Range("H46").Value = 1714
$ws.Range("B57").Value = 13953 ; $ws.Range("C57").Value = 397 ; $ws.Range("E57").Value = 5744 ; $ws.Range("G57").Value = 14 ; $ws.Range("H57").Value = 464
$ws.Range("B63").Value = 12238 ; $ws.Range("C63").Value = 471 ; $ws.Range("D63").Value = 6516 ; $ws.Range("E63").Value = 5574 ; $ws.Range("G63").Value = 5 ; $ws.Range("H63").Value = 148
$ws.Range("B80").Value = 5005 ; $ws.Range("C80").Value = 185 ; $ws.Range("D80").Value = 1904 ; $ws.Range("E80").Value = 2868 ; $ws.Range("G80").Value = 11 ; $ws.Range("H80").Value = 233
$ws.Range("B81").Value = 4980 ; $ws.Range("C81").Value = 64 ; $ws.Range("D81").Value = 24 ; $ws.Range("E81").Value = 4869 ; $ws.Range("G81").Value = 3 ; $ws.Range("H81").Value = 87
$ws.Range("B82").Value = 4904 ; $ws.Range("D82").Value = 3522 ; $ws.Range("E82").Value = 1355 ; $ws.Range("H82").Value = 27
$ws.Range("E84").Value = 2771 ; $ws.Range("G84").Value = 2 ; $ws.Range("H84").Value = 121
$ws.Range("B86").Value = 4469 ; $ws.Range("C86").Value = 399 ; $ws.Range("D86").Value = 1029 ; $ws.Range("E86").Value = 3368 ; $ws.Range("H86").Value = 72
$ws.Range("B87").Value = 4428 ; $ws.Range("D87").Value = 1750 ; $ws.Range("E87").Value = 2644 ; $ws.Range("H87").Value = 34
$ws.Range("B88").Value = 4099 ; $ws.Range("C88").Value = 0 ; $ws.Range("D88").Value = 3944 ; $ws.Range("E88").Value = 45 ; $ws.Range("G88").Value = 0 ; $ws.Range("H88").Value = 110
$ws.Range("B89").Value = 4086 ; $ws.Range("C89").Value = 5 ; $ws.Range("D89").Value = 2585 ; $ws.Range("E89").Value = 931 ; $ws.Range("G89").Value = 2 ; $ws.Range("H89").Value = 570
$ws.Range("B93").Value = 3256 ; $ws.Range("C93").Value = 19 ; $ws.Range("E93").Value = 1692 ; $ws.Range("G93").Value = 1 ; $ws.Range("H93").Value = 190
$ws.Range("B102").Value = 2163 ; $ws.Range("C102").Value = 194 ; $ws.Range("D102").Value = 890 ; $ws.Range("E102").Value = 1268 ; $ws.Range("H102").Value = 5
$ws.Range("B103").Value = 2150 ; $ws.Range("D103").Value = 1769 ; $ws.Range("E103").Value = 373 ; $ws.Range("H103").Value = 8
$ws.Range("B104").Value = 2058 ; $ws.Range("C104").Value = 0 ; $ws.Range("D104").Value = 982 ; $ws.Range("E104").Value = 1064 ; $ws.Range("H104").Value = 12
$ws.Range("B105").Value = 1981 ; $ws.Range("C105").Value = 2 ; $ws.Range("D105").Value = 1758 ; $ws.Range("E105").Value = 154 ; $ws.Range("H105").Value = 69
$ws.Range("E124").Value = 47 ; $ws.Range("G124").Value = 1 ; $ws.Range("H124").Value = 5
$ws.Range("B133").Value = 863 ; $ws.Range("C133").Value = 15 ; $ws.Range("D133").Value = 377 ; $ws.Range("E133").Value = 478 ; $ws.Range("H133").Value = 8
$ws.Range("B134").Value = 858 ; $ws.Range("D134").Value = 742 ; $ws.Range("E134").Value = 42 ; $ws.Range("H134").Value = 74
$ws.Range("B135").Value = 855 ; $ws.Range("D135").Value = 792 ; $ws.Range("E135").Value = 11 ; $ws.Range("H135").Value = 52
$ws.Range("B136").Value = 853 ; $ws.Range("D136").Value = 814 ; $ws.Range("E136").Value = 15 ; $ws.Range("H136").Value = 24
$ws.Range("B149").Value = 601 ; $ws.Range("C149").Value = 20 ; $ws.Range("D149").Value = 254 ; $ws.Range("E149").Value = 314
$ws.Range("D210").Value = 10 ; $ws.Range("H210").Value = 1
$ws.Range("D211").Value = 11 ; $ws.Range("H211").Value = 0
$ws.Range("D213").Value = 7 ; $ws.Range("H213").Value = 1
$ws.Range("D214").Value = 8 ; $ws.Range("H214").Value = 0
